$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2021" column (O) to the waste-disposal table, mirroring ---
# --- the formatting already used for the equivalent row in column N.   ---

# Row 2: blank separator cell, same border/style as its neighbor N2.
$ws.Range("N2").Copy() | Out-Null
$ws.Range("O2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 3: year header "2021", same style as N3.
$ws.Range("N3").Copy() | Out-Null
$ws.Range("O3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("O3").Value = 2021

# Row 4: per-person waste figure for 2021, computed from O5/O6, same style as N4.
$ws.Range("N4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("O4").Formula = "=O5/O6*1000"

# Row 5: total removed waste (thousand tons) for 2021 - style matches the
# plain (unformatted-number) style already used by B5/C5 in this row.
$ws.Range("B5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("O5").Value = 1229.5999999999999

# Row 6: average annual population for 2021, same style as N6.
$ws.Range("N6").Copy() | Out-Null
$ws.Range("O6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("O6").Value = 6436.9

$excel.CutCopyMode = $false

# Move / record the active selection as in the authored workbook.
$ws.Range("P16").Select() | Out-Null
